$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Displacement Table" sheet - update Displacement (tonnes) column B
#    for rows 2..26 (Draft 2m .. 14m)
# ---------------------------------------------------------------------
$wsDisp = $wb.Worksheets.Item("Displacement Table")

$newDisplacement = @(801, 1400, 2208, 3247, 4533, 6085, 7919, 10050, 12492, 15260, 18366, 21823, 25644, 29841, 34425, 39407, 44799, 50610, 56852, 63534, 70667, 78260, 86322, 94863, 103892)

for ($i = 0; $i -lt $newDisplacement.Length; $i++) {
    $rowNum = $i + 2
    $wsDisp.Cells.Item($rowNum, 2).Value2 = $newDisplacement[$i]
}

# ---------------------------------------------------------------------
# 2) "KN Curves" sheet - replace the displacement/KN data table
#    (rows 2..41, columns A..L). The table grows from 32 data rows
#    (rows 2..33) to 40 data rows (rows 2..41).
# ---------------------------------------------------------------------
$wsKN = $wb.Worksheets.Item("KN Curves")

$kn = @(
    @(1000, 0.51, 0.98, 1.58, 2.3, 3.14, 4.08, 5.11, 6.22, 7.38, 8.57, 9.789999999999999),
    @(3795, 1.13, 2.17, 3.51, 5.12, 6.99, 9.09, 11.38, 13.84, 16.42, 19.08, 21.79),
    @(6590, 1.58, 3.02, 4.88, 7.13, 9.73, 12.65, 15.85, 19.27, 22.86, 26.58, 30.34),
    @(9385, 1.95, 3.74, 6.04, 8.81, 12.03, 15.64, 19.59, 23.82, 28.27, 32.86, 37.51),
    @(12179, 2.28, 4.37, 7.06, 10.31, 14.07, 18.29, 22.91, 27.86, 33.05, 38.42, 43.87),
    @(14974, 2.58, 4.95, 7.99, 11.67, 15.93, 20.71, 25.93, 31.53, 37.41, 43.49, 49.65),
    @(17769, 2.86, 5.48, 8.85, 12.93, 17.65, 22.95, 28.74, 34.94, 41.46, 48.19, 55.02),
    @(20564, 3.12, 5.98, 9.66, 14.11, 19.27, 25.05, 31.37, 38.14, 45.26, 52.6, 60.06),
    @(23359, 3.37, 6.46, 10.43, 15.23, 20.8, 27.04, 33.86, 41.17, 48.85, 56.78, 64.84),
    @(26154, 3.61, 6.91, 11.16, 16.3, 22.26, 28.93, 36.24, 44.06, 52.28, 60.77, 69.38),
    @(28949, 3.83, 7.35, 11.87, 17.33, 23.65, 30.75, 38.52, 46.83, 55.57, 64.58, 73.73999999999999),
    @(31744, 4.05, 7.76, 12.54, 18.31, 25, 32.5, 40.71, 49.49, 58.73, 68.26000000000001, 77.94),
    @(34538, 4.26, 8.17, 13.19, 19.26, 26.3, 34.19, 42.82, 52.06, 61.77, 71.8, 81.98),
    @(37333, 4.47, 8.56, 13.82, 20.18, 27.55, 35.82, 44.87, 54.55, 64.73, 75.23, 85.90000000000001),
    @(40128, 4.66, 8.94, 14.43, 21.08, 28.77, 37.41, 46.85, 56.97, 67.59, 78.56, 89.7),
    @(42923, 4.86, 9.300000000000001, 15.03, 21.95, 29.96, 38.95, 48.78, 59.32, 70.38, 81.8, 93.40000000000001),
    @(45718, 5.04, 9.66, 15.61, 22.79, 31.11, 40.45, 50.67, 61.6, 73.09, 84.95999999999999, 97.01000000000001),
    @(48513, 5.23, 10.01, 16.17, 23.62, 32.24, 41.92, 52.5, 63.84, 75.73999999999999, 88.04000000000001, 100.52),
    @(51308, 5.4, 10.36, 16.73, 24.43, 33.34, 43.35, 54.3, 66.02, 78.33, 91.05, 103.96),
    @(54103, 5.58, 10.69, 17.27, 25.22, 34.42, 44.75, 56.05, 68.15000000000001, 80.87, 93.98999999999999, 107.32),
    @(56897, 5.75, 11.02, 17.8, 25.99, 35.48, 46.13, 57.77, 70.23999999999999, 83.34999999999999, 96.87, 110.61),
    @(59692, 5.92, 11.34, 18.32, 26.75, 36.51, 47.47, 59.46, 72.29000000000001, 85.78, 99.7, 113.84),
    @(62487, 6.08, 11.66, 18.83, 27.49, 37.53, 48.79, 61.11, 74.31, 88.17, 102.48, 117.01),
    @(65282, 6.24, 11.97, 19.33, 28.22, 38.53, 50.09, 62.74, 76.28, 90.51000000000001, 105.2, 120.12),
    @(68077, 6.4, 12.27, 19.82, 28.94, 39.51, 51.37, 64.34, 78.23, 92.81999999999999, 107.88, 123.18),
    @(70872, 6.56, 12.57, 20.3, 29.65, 40.48, 52.62, 65.91, 80.14, 95.09, 110.52, 126.19),
    @(73667, 6.71, 12.87, 20.78, 30.35, 41.43, 53.86, 67.45999999999999, 82.02, 97.31999999999999, 113.11, 129.15),
    @(76462, 6.87, 13.16, 21.25, 31.03, 42.36, 55.08, 68.98, 83.87, 99.52, 115.67, 132.07),
    @(79256, 7.02, 13.44, 21.71, 31.71, 43.28, 56.27, 70.48, 85.7, 101.68, 118.19, 134.95),
    @(82051, 7.16, 13.73, 22.17, 32.37, 44.19, 57.46, 71.95999999999999, 87.5, 103.82, 120.67, 137.78),
    @(84846, 7.31, 14, 22.62, 33.03, 45.09, 58.62, 73.43000000000001, 89.28, 105.93, 123.12, 140.58),
    @(87641, 7.45, 14.28, 23.06, 33.68, 45.98, 59.77, 74.87, 91.03, 108.01, 125.54, 143.34),
    @(90436, 7.59, 14.55, 23.5, 34.32, 46.85, 60.91, 76.29000000000001, 92.76000000000001, 110.06, 127.93, 146.07),
    @(93231, 7.73, 14.82, 23.93, 34.95, 47.71, 62.03, 77.7, 94.47, 112.09, 130.28, 148.76),
    @(96026, 7.87, 15.08, 24.36, 35.58, 48.57, 63.14, 79.09, 96.16, 114.09, 132.61, 151.42),
    @(98821, 8.01, 15.35, 24.79, 36.2, 49.41, 64.23999999999999, 80.45999999999999, 97.83, 116.08, 134.92, 154.05),
    @(101615, 8.140000000000001, 15.6, 25.2, 36.81, 50.24, 65.31999999999999, 81.81999999999999, 99.48, 118.03, 137.19, 156.65),
    @(104410, 8.279999999999999, 15.86, 25.62, 37.41, 51.07, 66.39, 83.16, 101.11, 119.97, 139.44, 159.22),
    @(107205, 8.41, 16.11, 26.03, 38.01, 51.88, 67.45999999999999, 84.48999999999999, 102.73, 121.89, 141.67, 161.76),
    @(110000, 8.539999999999999, 16.37, 26.43, 38.6, 52.69, 68.51000000000001, 85.8, 104.33, 123.78, 143.88, 164.28)
)

for ($i = 0; $i -lt $kn.Length; $i++) {
    $rowNum = $i + 2
    $rowData = $kn[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $wsKN.Cells.Item($rowNum, $j + 1).Value2 = $rowData[$j]
    }
}
